# Reassign trial rows 2-41 of the scene-categorization sheet using a fixed
# row permutation (new content for each destination row is taken verbatim
# from a source row, per the commit: "make only 20 different versions and
# duplicate many times for 1000 subjects"). Column F (trial_total) is
# recomputed as destination_row + 404. Columns A-E, G, J are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("H","I","K","L","M","N","O","P","Q","R","S","T","U","V")

# destination row -> source row (both 2..41), derived from the target diff.
$mapping = @(
    @(2, 37),
    @(3, 16),
    @(4, 25),
    @(5, 33),
    @(6, 26),
    @(7, 32),
    @(8, 12),
    @(9, 2),
    @(10, 22),
    @(11, 10),
    @(12, 7),
    @(13, 41),
    @(14, 36),
    @(15, 6),
    @(16, 15),
    @(17, 38),
    @(18, 18),
    @(19, 5),
    @(20, 29),
    @(21, 39),
    @(22, 17),
    @(23, 4),
    @(24, 40),
    @(25, 24),
    @(26, 28),
    @(27, 21),
    @(28, 27),
    @(29, 34),
    @(30, 20),
    @(31, 9),
    @(32, 13),
    @(33, 30),
    @(34, 23),
    @(35, 3),
    @(36, 19),
    @(37, 14),
    @(38, 35),
    @(39, 11),
    @(40, 8),
    @(41, 31)
)

# Snapshot every source row's relevant cell values before writing anything,
# since the mapping is a permutation and rows are both sources and
# destinations (writing in place could clobber data still needed later).
$snapshot = @{}
for ($r = 2; $r -le 41; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range($c + $r).Value2
    }
    $snapshot[$r] = $rowData
}

foreach ($pair in $mapping) {
    $destRow = $pair[0]
    $srcRow = $pair[1]
    $data = $snapshot[$srcRow]

    foreach ($c in $cols) {
        $ws.Range($c + $destRow).Value2 = $data[$c]
    }

    $ws.Range("F" + $destRow).Value2 = $destRow + 404
}
